# Updates for 21 release
# - Back up the current "CV" sheet to a new sheet named "bak" (full copy of
#   the CV list as it stood before this edit).
# - Remove a handful of retired template entries from the live "CV" sheet,
#   leaving blank rows where they used to be (the list shifts up).

$wb = $excel.ActiveWorkbook
$cv = $wb.Worksheets.Item("CV")

# 1) Duplicate CV -> "bak" (placed after the last sheet) before trimming it.
$cv.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$bak = $wb.Worksheets.Item($wb.Worksheets.Count)
$bak.Name = "bak"

# Leave the backup sheet's own view roughly where Excel would park it after
# a copy (top of the list, nothing special selected) before returning focus
# to CV below.
$bak.Activate()
$bak.Range("A5").Select()

# 2) Remove the retired templates from the live CV sheet. Deleting the
# entire row shifts everything below it up and appends a blank row at the
# bottom, which is exactly the before/after shape we need.
$namesToRemove = @(
    "ConvertToTiff",
    "ExtractHtml",
    "htmlToPdf stream test",
    "s3keys",
    "XMLTransform"
)

foreach ($name in $namesToRemove) {
    $found = $cv.Columns.Item(1).Find($name)
    if ($found -ne $null) {
        $found.EntireRow.Delete()
    }
}

# 3) Restore CV as the active sheet/selection, matching the saved view.
$cv.Activate()
$cv.Range("A11").Select()
